# Populate the previously-empty (all-zero / placeholder "-") quarterly
# income-statement figures on the "Overview" sheet with the real reported
# values (rows 11-27, columns D:M = the 10 fiscal quarters).
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Overview")

# Row 11: فروش (Sales)
$ws.Range("D11").Value = 560992
$ws.Range("E11").Value = 750320
$ws.Range("F11").Value = 722753
$ws.Range("G11").Value = 807734
$ws.Range("H11").Value = 1108503
$ws.Range("I11").Value = 1078860
$ws.Range("J11").Value = 1483266
$ws.Range("K11").Value = 1903632
$ws.Range("L11").Value = 3112367
$ws.Range("M11").Value = 3188078

# Row 12: بهای تمام شده کالای فروش رفته (Cost of goods sold)
$ws.Range("D12").Value = -389089
$ws.Range("E12").Value = -482632
$ws.Range("F12").Value = -452361
$ws.Range("G12").Value = -457104
$ws.Range("H12").Value = -708691
$ws.Range("I12").Value = -619254
$ws.Range("J12").Value = -761642
$ws.Range("K12").Value = -1136800
$ws.Range("L12").Value = -1657837
$ws.Range("M12").Value = -1657608

# Row 13: سود (زیان) ناخالص (Gross profit)
$ws.Range("D13").Value = 171903
$ws.Range("E13").Value = 267688
$ws.Range("F13").Value = 270392
$ws.Range("G13").Value = 350630
$ws.Range("H13").Value = 399812
$ws.Range("I13").Value = 459606
$ws.Range("J13").Value = 721624
$ws.Range("K13").Value = 766832
$ws.Range("L13").Value = 1454530
$ws.Range("M13").Value = 1530470

# Row 14: هزینه های عمومی, اداری و تشکیلاتی (General & admin expenses)
$ws.Range("D14").Value = -18924
$ws.Range("E14").Value = -16690
$ws.Range("F14").Value = -23073
$ws.Range("G14").Value = -22232
$ws.Range("H14").Value = -20908
$ws.Range("I14").Value = -21366
$ws.Range("J14").Value = -55673
$ws.Range("K14").Value = -22360
$ws.Range("L14").Value = -83991
$ws.Range("M14").Value = -214651

# Row 15: هزینه کاهش ارزش دریافتنی‌ها (was the placeholder text "-" in
# every quarter; now numeric 0, so the now-unused "-" shared string gets
# dropped automatically)
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 0

# Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی (Other operating income/expense, net)
$ws.Range("D16").Value = -15752
$ws.Range("E16").Value = -20614
$ws.Range("F16").Value = -29004
$ws.Range("G16").Value = -13307
$ws.Range("H16").Value = 11491
$ws.Range("I16").Value = 7288
$ws.Range("J16").Value = -64132
$ws.Range("K16").Value = -12556
$ws.Range("L16").Value = 20545
$ws.Range("M16").Value = -46701

# Row 17: سود (زیان) عملیاتی (Operating profit/loss)
$ws.Range("D17").Value = 137227
$ws.Range("E17").Value = 230384
$ws.Range("F17").Value = 218315
$ws.Range("G17").Value = 315091
$ws.Range("H17").Value = 390395
$ws.Range("I17").Value = 445528
$ws.Range("J17").Value = 601819
$ws.Range("K17").Value = 731916
$ws.Range("L17").Value = 1391084
$ws.Range("M17").Value = 1269118

# Row 18: هزینه های مالی (Financial expenses)
$ws.Range("D18").Value = -9997
$ws.Range("E18").Value = -9256
$ws.Range("F18").Value = 3011
$ws.Range("G18").Value = -5655
$ws.Range("H18").Value = -9807
$ws.Range("I18").Value = -10592
$ws.Range("J18").Value = -27778
$ws.Range("K18").Value = -16279
$ws.Range("L18").Value = -94374
$ws.Range("M18").Value = -117832

# Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی (Other non-operating income/expense, net)
$ws.Range("D19").Value = 202017
$ws.Range("E19").Value = 3926
$ws.Range("F19").Value = -41280
$ws.Range("G19").Value = 14163
$ws.Range("H19").Value = 838
$ws.Range("I19").Value = -15001
$ws.Range("J19").Value = 142564
$ws.Range("K19").Value = 19757
$ws.Range("L19").Value = 13826
$ws.Range("M19").Value = 60169

# Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات (Pre-tax profit from continuing ops)
$ws.Range("D20").Value = 329247
$ws.Range("E20").Value = 225054
$ws.Range("F20").Value = 180046
$ws.Range("G20").Value = 323599
$ws.Range("H20").Value = 381426
$ws.Range("I20").Value = 419935
$ws.Range("J20").Value = 716605
$ws.Range("K20").Value = 735394
$ws.Range("L20").Value = 1310536
$ws.Range("M20").Value = 1211455

# Row 21: مالیات (Tax) -- note G21 was also a placeholder "-", now 0
$ws.Range("D21").Value = -71496
$ws.Range("E21").Value = 71496
$ws.Range("F21").Value = -85675
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = -110831
$ws.Range("I21").Value = -52878
$ws.Range("J21").Value = 147193
$ws.Range("K21").Value = -73539
$ws.Range("L21").Value = -2095
$ws.Range("M21").Value = -36990

# Row 22: سود (زیان) خالص عملیات در حال تداوم (Net profit from continuing ops)
$ws.Range("D22").Value = 257751
$ws.Range("E22").Value = 296550
$ws.Range("F22").Value = 94371
$ws.Range("G22").Value = 323599
$ws.Range("H22").Value = 270595
$ws.Range("I22").Value = 367057
$ws.Range("J22").Value = 863798
$ws.Range("K22").Value = 661855
$ws.Range("L22").Value = 1308441
$ws.Range("M22").Value = 1174465

# Row 23: سود (زیان) عملیات متوقف شده پس از اثر مالیاتی (Discontinued ops) -- all placeholder "-" -> 0
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 0

# Row 24: سود (زیان) خالص (Net profit/loss)
$ws.Range("D24").Value = 257751
$ws.Range("E24").Value = 296550
$ws.Range("F24").Value = 94371
$ws.Range("G24").Value = 323599
$ws.Range("H24").Value = 270595
$ws.Range("I24").Value = 367057
$ws.Range("J24").Value = 863798
$ws.Range("K24").Value = 661855
$ws.Range("L24").Value = 1308441
$ws.Range("M24").Value = 1174465

# Row 25: سود هر سهم پس از کسر مالیات (EPS after tax)
$ws.Range("D25").Value = 175
$ws.Range("E25").Value = 202
$ws.Range("F25").Value = 64
$ws.Range("G25").Value = 63
$ws.Range("H25").Value = 53
$ws.Range("I25").Value = 71
$ws.Range("J25").Value = 168
$ws.Range("K25").Value = 129
$ws.Range("L25").Value = 255
$ws.Range("M25").Value = 228

# Row 26: سرمایه (Capital)
$ws.Range("D26").Value = 1468924
$ws.Range("E26").Value = 1468924
$ws.Range("F26").Value = 1468924
$ws.Range("G26").Value = 5141234
$ws.Range("H26").Value = 5141234
$ws.Range("I26").Value = 5141234
$ws.Range("J26").Value = 5141234
$ws.Range("K26").Value = 5141234
$ws.Range("L26").Value = 5141234
$ws.Range("M26").Value = 5141234

# Row 27: سود هر سهم بر اساس آخرین سرمایه (EPS based on latest capital)
$ws.Range("D27").Value = 50
$ws.Range("E27").Value = 58
$ws.Range("F27").Value = 18
$ws.Range("G27").Value = 63
$ws.Range("H27").Value = 53
$ws.Range("I27").Value = 71
$ws.Range("J27").Value = 168
$ws.Range("K27").Value = 129
$ws.Range("L27").Value = 255
$ws.Range("M27").Value = 228
